# v0.7.5b: Initial Ahimdoor map outlines
# Insert a new "Gazebo" / building feature row into the Features sheet,
# and bump a handful of flammability/burn-time ("E" column) values on the
# nearby "large tree" rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Features")

# Insert a new blank row at row 87 (pushes the old row 88.. block down by one,
# inheriting formatting from the row above exactly like an interactive
# Excel "Insert Row" would).
$ws.Rows.Item(87).Insert()

# Populate the new row 87 with the Gazebo/building feature entry.
$ws.Cells.Item(87, 2).Value = 253
$ws.Cells.Item(87, 3).Value = "Gazebo"
$ws.Cells.Item(87, 5).Value = 7
$ws.Cells.Item(87, 7).Value = "building"

# Column G has no sheet-level default style, so the quote-prefixed format
# that "Insert Row" copied down from row 86 gets dropped as soon as the
# value is written; restore it by re-pasting the format from the row above.
$ws.Cells.Item(86, 7).Copy()
$ws.Cells.Item(87, 7).PasteSpecial(-4122)
$excel.CutCopyMode = 0

# A handful of existing rows (now shifted down by one) had their "E" column
# burn-time values edited as part of this pass.
$ws.Cells.Item(135, 5).Value = 30
$ws.Cells.Item(141, 5).Value = 30
$ws.Cells.Item(142, 5).Value = 30
$ws.Cells.Item(143, 5).Value = 30
$ws.Cells.Item(144, 5).Value = 30
$ws.Cells.Item(145, 5).Value = 40
$ws.Cells.Item(146, 5).Value = 30

# Restore the view state described by the commit (scrolled back up to the
# newly-added row, which is now the active selection).
$ws.Application.ActiveWindow.ScrollRow = 73
$ws.Range("H87").Select()
